$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.586.16"
$ws.Range("E2").Value = "  -1.21%  "
$ws.Range("D3").Value = "2.063.58"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.78"
$ws.Range("E5").Value = "  -2.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.676"
$ws.Range("E6").Value = "  +2.04%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.59"
$ws.Range("E8").Value = "  -6.59%  "
$ws.Range("E9").Value = "  -1.72%  "
$ws.Range("E10").Value = "  -5.79%  "
$ws.Range("E11").Value = "  -3.56%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.107"
$ws.Range("E12").Value = "  -0.96%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.889"
$ws.Range("E13").Value = "  -0.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.41"
$ws.Range("E14").Value = "  -9.56%  "
$ws.Range("D15").Value = "2.366.46"
$ws.Range("E15").Value = "  +1.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.39"
$ws.Range("E16").Value = "  -5.02%  "
$ws.Range("D17").Value = "2.070.79"
$ws.Range("E17").Value = "  +1.54%  "
$ws.Range("D18").Value = "36.497.18"
$ws.Range("E18").Value = "  -1.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.34"
$ws.Range("E19").Value = "  -12.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.61"
$ws.Range("E20").Value = "  -4.16%  "
$ws.Range("D21").Value = "0.0₃0862"
$ws.Range("E21").Value = "  -2.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.26"
$ws.Range("E22").Value = "  -2.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.54"
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("E25").Value = "  -4.96%  "
$ws.Range("E26").Value = "  -2.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.24"
$ws.Range("E27").Value = "  -3.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.08"
$ws.Range("E28").Value = "  -4.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.28"
$ws.Range("E29").Value = "  +0.99%  "
$ws.Range("E30").Value = "  -0.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.06"
$ws.Range("E31").Value = "  -0.89%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.14"
$ws.Range("E32").Value = "  -3.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.57"
$ws.Range("E33").Value = "  -1.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0595"
$ws.Range("E34").Value = "  -4.04%  "
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.27"
$ws.Range("E36").Value = "  +0.79%  "
$ws.Range("E37").Value = "  -1.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0813"
$ws.Range("E38").Value = "  -6.89%  "
$ws.Range("E39").Value = "  -6.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.86"
$ws.Range("E40").Value = "  -5.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.89"
$ws.Range("E41").Value = "  -6.13%  "
$ws.Range("E42").Value = "  -3.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.12"
$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0927"
$ws.Range("E44").Value = "  -6.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.57"
$ws.Range("E45").Value = "  -4.38%  "
$ws.Range("D46").Value = "1.385.54"
$ws.Range("E46").Value = "  +8.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "15.57"
$ws.Range("E47").Value = "  -9.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.29"
$ws.Range("E48").Value = "  +7.33%  "
$ws.Range("E49").Value = "  -3.24%  "
$ws.Range("D51").Value = "2.252.54"
$ws.Range("E51").Value = "  +1.33%  "
